# Berechnungen.xlsx edit script
# - change F_CPU from 14745600 to 16000000
# - rename Tabelle1 -> "SPI & Timer"
# - add Timer calc block (rows 7-13) on sheet1
# - add new worksheet "Widerständen" with resistor/potentiometer calc table

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: "SPI & Timer" ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "SPI & Timer"

# F_CPU changed to 16 MHz
$ws1.Range("B1").Value = 16000000

# Timer prescaler block
$ws1.Range("A7").Value = "Timer prescaler"
$ws1.Range("B7").Value = 1024

$ws1.Range("A8").Value = "Timer clock"
$ws1.Range("B8").Formula = "=B1/B7"

$ws1.Range("A9").Value = "Timer periode"
$ws1.Range("B9").Formula = "=1/B8"

$ws1.Range("A10").Value = "Timer periode, ms"
$ws1.Range("B10").Formula = "=B9*1000"

$ws1.Range("A12").Value = "Gesuchte Periode, ms"
$ws1.Range("B12").Value = 1.7

$ws1.Range("A13").Value = "Counter"
$ws1.Range("B13").Formula = "=B12/B10"
$ws1.Range("D13").Value = "bei 16 MHz = 18"
$ws1.Range("F13").Value = "bei 14,745600 = 24"

$ws1.Columns.Item(1).AutoFit()

$ws1.Range("B10").Select() | Out-Null

# ---- Sheet 2: "Widerständen" ----
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Widerständen"

$ws2.Range("A1").Value = "Wert"
$ws2.Range("B1").Value = "POT1"
$ws2.Range("C1").Value = "POT2"
$ws2.Range("D1").Value = "Berechnet"
$ws2.Range("E1").Value = "Wirklichkeit"

$ws2.Range("A2").Value = "2k2"
$ws2.Range("B2").Value = 250
$ws2.Range("C2").Value = 251
$ws2.Range("D2").Formula = "=(((50000 * (256-B2))/256)+52)+(((50000 * (256-C2))/256)+52)"
$ws2.Range("E2").Value = 2500

$ws2.Range("A3").Value = "4k4"
$ws2.Range("B3").Value = 240
$ws2.Range("C3").Value = 250
$ws2.Range("D3").Formula = "=(((50000 * (256-B3))/256)+52)+(((50000 * (256-C3))/256)+52)"
$ws2.Range("E3").Value = 4667

$ws2.Range("A4").Value = "8k8"
$ws2.Range("B4").Value = 237
$ws2.Range("C4").Value = 230
$ws2.Range("D4").Formula = "=(((50000 * (256-B4))/256)+52)+(((50000 * (256-C4))/256)+52)"
$ws2.Range("E4").Value = 9200

$ws2.Range("A5").Value = "12k1"
$ws2.Range("B5").Value = 226
$ws2.Range("C5").Value = 225
$ws2.Range("D5").Formula = "=(((50000 * (256-B5))/256)+52)+(((50000 * (256-C5))/256)+52)"
$ws2.Range("E5").Value = 12370

$ws2.Range("A6").Value = "16k8"
$ws2.Range("B6").Value = 214
$ws2.Range("C6").Value = 214
$ws2.Range("D6").Formula = "=(((50000 * (256-B6))/256)+52)+(((50000 * (256-C6))/256)+52)"
$ws2.Range("E6").Value = 16910

$ws2.Range("A7").Value = "23k6"
$ws2.Range("B7").Value = 196
$ws2.Range("C7").Value = 196
$ws2.Range("D7").Formula = "=(((50000 * (256-B7))/256)+52)+(((50000 * (256-C7))/256)+52)"
$ws2.Range("E7").Value = 24000
$ws2.Range("G7").Formula = "=50000/256"

$ws2.Range("A8").Value = "6k6"
$ws2.Range("B8").Value = 240
$ws2.Range("C8").Value = 240
$ws2.Range("D8").Formula = "=(((50000 * (256-B8))/256)+52)+(((50000 * (256-C8))/256)+52)"
$ws2.Range("E8").Value = 6650

$ws2.Range("A9").Value = "33k6"
$ws2.Range("B9").Value = 162
$ws2.Range("C9").Value = 182
$ws2.Range("D9").Formula = "=(((50000 * (256-B9))/256)+52)+(((50000 * (256-C9))/256)+52)"
$ws2.Range("E9").Value = 33500

$ws2.Range("A10").Value = "48k6"
$ws2.Range("B10").Value = 135
$ws2.Range("C10").Value = 135
$ws2.Range("D10").Formula = "=(((50000 * (256-B10))/256)+52)+(((50000 * (256-C10))/256)+52)"
$ws2.Range("E10").Value = 48100

$ws2.Columns.Item(4).AutoFit()
$ws2.Columns.Item(5).AutoFit()

# Sheet2 uses ~2cm top/bottom page margins instead of the 0.75in default
$ws2.PageSetup.TopMargin = 56.69291339999999
$ws2.PageSetup.BottomMargin = 56.69291339999999

$ws2.Range("C10").Select() | Out-Null

# Re-activate sheet 1 so it is the selected tab on save (matches original file)
$ws1.Activate()
$ws1.Range("B10").Select() | Out-Null

$excel.Windows.Item(1).Top = 1905
$excel.Windows.Item(1).Left = 240
